# NumerosSelecionados.xlsx — "ESTÉTICA E BOTÕES OK"
#   1) row 23: idPagamento (col D) was blank, now filled in
#   2) a new row 29 is appended: another "Vitor Ito" submission
#
# Columns C (Telefone) and D (idPagamento) hold long numeric-looking ids
# that must stay TEXT (like every other row in this sheet), so each one is
# pre-formatted as Text ("@") before the value is assigned — otherwise
# Excel would auto-convert the numeric string into a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: fill in the previously-blank idPagamento ---------------------
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "76955304908"

# --- New row 29: another "Vitor Ito" entry, same shape as rows 9-28 -------
$ws.Range("A29").Value = "Vitor Ito"
$ws.Range("B29").Value = 1578424633
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "11966554433"
# idPagamento (D29) is left blank for this entry, same as rows 24-28
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 4
$ws.Range("I29").Value = 5
$ws.Range("J29").Value = 6
$ws.Range("K29").Value = 7
$ws.Range("L29").Value = 8
$ws.Range("M29").Value = 9
$ws.Range("N29").Value = 10
$ws.Range("O29").Value = "Não"
